$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the greeting text for rule R10 (cell E8) from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Select cell E8 so it becomes the active cell / selection in the sheet view
$ws.Range("E8").Select()
